$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: update K14 (value only, row already spans to column K)
$ws.Range("K14").Value = -0.2254024683979639

# Row 15: update J15, add new K15
$ws.Range("J15").Value = -0.1253231084953424
$ws.Range("K15").Value = -0.3352267436446591

# Row 16: update I16, add new J16
$ws.Range("I16").Value = 0.2284633975843539
$ws.Range("J16").Value = 0.01855976243503714

# Row 17: update H17, add new I17
$ws.Range("H17").Value = 0.08028600715190851
$ws.Range("I17").Value = -0.1296176279974082

# Row 18: update G18, add new H18
$ws.Range("G18").Value = -0.07715998185224648
$ws.Range("H18").Value = -0.2870636170015632

# Row 19: update F19, add new G19
$ws.Range("F19").Value = 0.4234994746738243
$ws.Range("G19").Value = 0.2135958395245076

# Row 20: update E20, add new F20
$ws.Range("E20").Value = 0.1431415941383551
$ws.Range("F20").Value = -0.06676204101096155

# Row 21: update D21, add new E21
$ws.Range("D21").Value = 0.3151164519833668
$ws.Range("E21").Value = 0.1052128168340501

# Row 22: update C22, add new D22
$ws.Range("C22").Value = 0.009253912237035311
$ws.Range("D22").Value = -0.2006497229122814

# Row 23: update B23, add new C23
$ws.Range("B23").Value = 0.6215838649243215
$ws.Range("C23").Value = 0.4116802297750048

# Row 24: add new B24
$ws.Range("B24").Value = -0.2766911554241067
